$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPP < 70kmph")

# 1) Shift the existing Rank column (C2:C321) up by 59 (1 -> 60, ..., 320 -> 379)
$oldRange = $ws.Range("C2:C321")
$oldVals = $oldRange.Value()
for ($i = 1; $i -le 320; $i++) {
    $oldVals[$i, 1] = $oldVals[$i, 1] + 59
}
$oldRange.Value = $oldVals

# 2) Append 59 new rows (322-380): A,B,D,E,F = 0 ; C = row-1 (321..379)
$newRange = $ws.Range("A322:F380")
$newVals = $newRange.Value()
for ($i = 1; $i -le 59; $i++) {
    $rowNum = 321 + $i
    $newVals[$i, 1] = 0
    $newVals[$i, 2] = 0
    $newVals[$i, 3] = $rowNum - 1
    $newVals[$i, 4] = 0
    $newVals[$i, 5] = 0
    $newVals[$i, 6] = 0
}
$newRange.Value = $newVals
